$d = $word.ActiveDocument

# The paragraph currently reads:  TTT class: summary TTT
# with the hidden "_GoBack" bookmark sitting right after "summary".
#
# Target reading (per diff):      TTT class: s atisfactionS [bookmark] ummary TTT
# i.e. "summary" -> "satisfactionSummary", split into three runs
# ("s", "atisfactionS", "ummary") with the bookmark re-anchored between
# "atisfactionS" and "ummary".

# Locate "summary" precisely (robust to any future shifts).
$findRange = $d.Content
$found = $findRange.Find.Execute("summary", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$summaryStart = $findRange.Start
$summaryEnd = $findRange.End

# Position right after the leading "s" of "summary".
$afterS = $summaryStart + 1

# Insert the new middle chunk "atisfactionS" right after that leading "s".
$insertPoint = $d.Range($afterS, $afterS)
$insertPoint.InsertAfter("atisfactionS")

# New layout (character offsets): ...class:[s][atisfactionS][ummary]TTT...
$afterAtisfactionS = $afterS + "atisfactionS".Length

# Force a run break between "class:" and "s" by toggling a formatting
# property on just the "s" character and back to its original value.
$rS = $d.Range($summaryStart, $afterS)
$rS.Bold = 1
$rS.Bold = 0

# Force a run break between "atisfactionS" and "ummary" the same way.
$rAtisfactionS = $d.Range($afterS, $afterAtisfactionS)
$rAtisfactionS.Bold = 1
$rAtisfactionS.Bold = 0

# Re-anchor the "_GoBack" bookmark so it sits between "atisfactionS" and "ummary".
$bookmark = $d.Bookmarks("_GoBack")
$bookmark.Delete()
$newBookmarkRange = $d.Range($afterAtisfactionS, $afterAtisfactionS)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)

Write-Output "Final paragraph text: $($d.Content.Text)"
